$wb = $excel.ActiveWorkbook

# --- Sheet: P_valores ---
$wsP = $wb.Worksheets.Item("P_valores")

$wsP.Range("C2").Value = 0.8612186940460373
$wsP.Range("D2").Value = 0.9925493497395685
$wsP.Range("E2").Value = 0.7920058109888606
$wsP.Range("F2").Value = 0.9738881191490687

$wsP.Range("B3").Value = 0.8612186940460373
$wsP.Range("D3").Value = 0.7515266661884532
$wsP.Range("E3").Value = 0.707065613347565
$wsP.Range("F3").Value = 0.8590567467227268

$wsP.Range("B4").Value = 0.9925493497395685
$wsP.Range("C4").Value = 0.7515266661884532
$wsP.Range("E4").Value = 0.7838537058857651
$wsP.Range("F4").Value = 0.9769650480470065

$wsP.Range("B5").Value = 0.7920058109888606
$wsP.Range("C5").Value = 0.707065613347565
$wsP.Range("D5").Value = 0.7838537058857651
$wsP.Range("F5").Value = 0.8007318003495649

$wsP.Range("B6").Value = 0.9738881191490687
$wsP.Range("C6").Value = 0.8590567467227268
$wsP.Range("D6").Value = 0.9769650480470065
$wsP.Range("E6").Value = 0.8007318003495649

# --- Sheet: Estadisticos_DM ---
$wsD = $wb.Worksheets.Item("Estadisticos_DM")

$wsD.Range("C2").Value = 0.1780703071381367
$wsD.Range("D2").Value = -0.009506263065901771
$wsD.Range("E2").Value = -0.2687978693654273
$wsD.Range("F2").Value = 0.03332214120984043

$wsD.Range("B3").Value = -0.1780703071381367
$wsD.Range("D3").Value = -0.3229221721493704
$wsD.Range("E3").Value = -0.3835609956233201
$wsD.Range("F3").Value = -0.1808766411807107

$wsD.Range("B4").Value = 0.009506263065901771
$wsD.Range("C4").Value = 0.3229221721493704
$wsD.Range("E4").Value = -0.2796253561128079
$wsD.Range("F4").Value = 0.02939428867276965

$wsD.Range("B5").Value = 0.2687978693654273
$wsD.Range("C5").Value = 0.3835609956233201
$wsD.Range("D5").Value = 0.2796253561128079
$wsD.Range("F5").Value = 0.2572452254893996

$wsD.Range("B6").Value = -0.03332214120984043
$wsD.Range("C6").Value = 0.1808766411807107
$wsD.Range("D6").Value = -0.02939428867276965
$wsD.Range("E6").Value = -0.2572452254893996
